# chinh sua yeu cau in ra header_footer
# Rename the header/footer marker rows from "header"/"footer" to
# "begin header"/"begin footer" on Sheet1 (the "end header"/"end footer"
# markers are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" keeps the existing quote-prefix cell style (s="1") intact
# while updating only the text content.
$ws.Range("A1").Value = "'----------------begin header-----------------"
$ws.Range("A6").Value = "'----------------begin footer-----------------"

# Match the new active-cell selection recorded in the sheet view.
$ws.Range("M5").Select()
